$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2040.2963
$ws.Range("I100").Value = 1588.5714
$ws.Range("J100").Value = 2526.7693
$ws.Range("K100").Value = 1588.5714
$ws.Range("L100").Value = 2526.7693
$ws.Range("M100").Value = -1047.5714
$ws.Range("N100").Value = -3608.7693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 94539.16
$ws.Range("I45").Value = 101832.914
$ws.Range("J45").Value = 7014
$ws.Range("K45").Value = 101832.914
$ws.Range("L45").Value = 7014
$ws.Range("M45").Value = -101455.914
$ws.Range("N45").Value = -7768

$ws.Range("H102").Value = 58850330
$ws.Range("I102").Value = 83335590
$ws.Range("J102").Value = 85686.2
$ws.Range("K102").Value = 83335590
$ws.Range("L102").Value = 85686.2
$ws.Range("M102").Value = -83333968
$ws.Range("N102").Value = -88930.2

$ws.Range("H110").Value = 2075
$ws.Range("I110").Value = 2075
$ws.Range("K110").Value = 2075
$ws.Range("M110").Value = -30

$ws.Range("H122").Value = 2422
$ws.Range("I122").Value = 2422
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7266
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -4816

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2103
$ws.Range("I20").Value = 1318
$ws.Range("J20").Value = 2574
$ws.Range("K20").Value = 1318
$ws.Range("L20").Value = 2574
$ws.Range("M20").Value = -1071
$ws.Range("N20").Value = -3068

$ws.Range("H86").Value = 2942.4783
$ws.Range("I86").Value = 2986.5334
$ws.Range("K86").Value = 2986.5334
$ws.Range("M86").Value = -1863.5334

$ws.Range("H89").Value = 2942.4783
$ws.Range("I89").Value = 2986.5334
$ws.Range("K89").Value = 14932.667
$ws.Range("M89").Value = -9316.666999999999

$ws.Range("H94").Value = 3001.6
$ws.Range("I94").Value = 3002
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 3002
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -2551
$ws.Range("N94").Value = -3902

$ws.Range("H99").Value = 2441.7073
$ws.Range("I99").Value = 2350.3572
$ws.Range("J99").Value = 2638.4614
$ws.Range("K99").Value = 2350.3572
$ws.Range("L99").Value = 2638.4614
$ws.Range("M99").Value = -852.3571999999999
$ws.Range("N99").Value = -5634.4614

$ws.Range("H107").Value = 1809.3334
$ws.Range("I107").Value = 1283.875
$ws.Range("J107").Value = 6013
$ws.Range("K107").Value = 1283.875
$ws.Range("L107").Value = 6013
$ws.Range("M107").Value = 636.125
$ws.Range("N107").Value = -9853

$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774

$ws.Range("H112").Value = 51468.668
$ws.Range("J112").Value = 51468.668
$ws.Range("L112").Value = 51468.668
$ws.Range("N112").Value = -54422.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 32000.125
$ws.Range("J120").Value = 32000.125
$ws.Range("L120").Value = 32000.125
$ws.Range("N120").Value = -39258.125

$ws.Range("H121").Value = 33319.75
$ws.Range("J121").Value = 33319.75
$ws.Range("L121").Value = 33319.75
$ws.Range("N121").Value = -35939.75

$ws.Range("H132").Value = 47479.066
$ws.Range("I132").Value = 960.44446
$ws.Range("J132").Value = 111889.46
$ws.Range("K132").Value = 2881.33338
$ws.Range("L132").Value = 335668.38
$ws.Range("M132").Value = -351.33338
$ws.Range("N132").Value = -340728.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 2075
$ws.Range("I99").Value = 905
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2715
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -469
$ws.Range("N99").Value = -19492

$ws.Range("H106").Value = 376785.56
$ws.Range("J106").Value = 6249.8335
$ws.Range("L106").Value = 18749.5005
$ws.Range("N106").Value = -20641.5005

$ws.Range("H131").Value = 2331.75
$ws.Range("I131").Value = 12911.125
$ws.Range("J131").Value = 1218.1316
$ws.Range("K131").Value = 38733.375
$ws.Range("L131").Value = 3654.3948
$ws.Range("M131").Value = -33693.375
$ws.Range("N131").Value = -13734.3948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5282.7144
$ws.Range("I70").Value = 5400.421
$ws.Range("J70").Value = 5034.222
$ws.Range("K70").Value = 5400.421
$ws.Range("L70").Value = 5034.222
$ws.Range("M70").Value = -5130.421
$ws.Range("N70").Value = -5574.222

$ws.Range("H73").Value = 5282.7144
$ws.Range("I73").Value = 5400.421
$ws.Range("J73").Value = 5034.222
$ws.Range("K73").Value = 5400.421
$ws.Range("L73").Value = 5034.222
$ws.Range("M73").Value = -4464.421
$ws.Range("N73").Value = -6906.222

$ws.Range("H80").Value = 6910.3887
$ws.Range("I80").Value = 5822.1113
$ws.Range("J80").Value = 7998.6665
$ws.Range("K80").Value = 5822.1113
$ws.Range("L80").Value = 7998.6665
$ws.Range("M80").Value = -4824.1113
$ws.Range("N80").Value = -9994.666499999999

$ws.Range("H83").Value = 6910.3887
$ws.Range("I83").Value = 5822.1113
$ws.Range("J83").Value = 7998.6665
$ws.Range("K83").Value = 29110.5565
$ws.Range("L83").Value = 39993.3325
$ws.Range("M83").Value = -24118.5565
$ws.Range("N83").Value = -49977.3325

$ws.Range("H97").Value = 3506.8635
$ws.Range("I97").Value = 2525.8823
$ws.Range("K97").Value = 2525.8823
$ws.Range("M97").Value = -2029.8823

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1307.1428
$ws.Range("I93").Value = 1680
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 1680
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = -432
$ws.Range("N93").Value = -3596

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1300243.4
$ws.Range("I122").Value = 1430117.8
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4290353.4
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4287903.4
$ws.Range("N122").Value = -9400

$ws.Range("H126").Value = 2451645.8
$ws.Range("I126").Value = 2674495.2
$ws.Range("J126").Value = 300
$ws.Range("K126").Value = 8023485.600000001
$ws.Range("L126").Value = 900
$ws.Range("M126").Value = -8021015.600000001
$ws.Range("N126").Value = -5840

Write-Host "Applied all updates"
